$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet held ~355,997 rows of raw spatial-inventory data
# (see the sortState range A2:G355997). This commit trims the fixture
# down to a single placeholder data row so the test file is small
# enough to live in the repo without Git LFS.
$ws.Range("A2").Value = 2010

# Move/record the active selection onto the new row, matching the
# saved cursor position.
$ws.Range("A2").Select()
